# Update crypto price/volume figures (refreshed data pull)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.193.53'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.33%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.098.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.90%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.49%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.96'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.21%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -0.57%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0841'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.36%  '

$ws.Range("E11").Value = '  +0.25%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.409.78'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.77%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '14.68'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("E15").Value = '  +5.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.777'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.094.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.61%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.175.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("E19").Value = '  +0.52%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.46%  '

$ws.Range("E21").Value = '  +1.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.43'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '

$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '170.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.69%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.99%  '

$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.90%  '

$ws.Range("E31").Value = '  -0.30%  '

$ws.Range("E32").Value = '  +9.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.69'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.31%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0606'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.38%  '

$ws.Range("E36").Value = '  +0.44%  '

$ws.Range("E37").Value = '  +5.39%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.547.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '100.14'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.83%  '

$ws.Range("E43").Value = '  +0.15%  '

$ws.Range("E44").Value = '  +0.81%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0906'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.39%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.66%  '

$ws.Range("E47").Value = '  +0.80%  '

$ws.Range("E48").Value = '  +1.45%  '

$ws.Range("E49").Value = '  +1.93%  '

$ws.Range("E50").Value = '  +1.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.295.46'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.86%  '
